$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (data_source_id) values from 82 to 1 for all data rows
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1

# Update the active selection to D10
$ws.Range("D10").Select()
